$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "ECs" target-cluster row (old row 2) is removed from the data; the
# remaining rows (old rows 3 and 4) shift up to rows 2 and 3. Deleting the
# entire row handles the row shift, the dimension shrink (A1:T4 -> A1:T3)
# and drops the now-unused "ECs" shared string automatically.
$ws.Rows(2).Delete()

# The raw expression columns carry over unchanged from the old rows, but the
# derived-specificity columns (O, P, S, T) are recomputed against the new
# (smaller) set of rows, so set their updated values explicitly.
$ws.Range("O2").Value = 0.8627666706462545
$ws.Range("P2").Value = 0.8627666706462543
$ws.Range("S2").Value = 0.8627666706462545
$ws.Range("T2").Value = 0.8627666706462543

$ws.Range("O3").Value = 0.1372333293537455
$ws.Range("P3").Value = 0.1372333293537455
$ws.Range("S3").Value = 0.1372333293537455
$ws.Range("T3").Value = 0.1372333293537455
